# Update "Demanda interna" workbook with the Agosto.2021 quarterly release.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BH: header "Agosto.2021" -------------------------------
# Give it the same (bold / bordered / centered) style as the other headers
# by copying the format from the previous header cell (BG1).
$ws.Range("BG1").Copy() | Out-Null
$ws.Range("BH1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("BH1").Value2 = "Agosto.2021"

# Carry forward the previous release's values (column BG) into the new
# column BH for every existing data row except row 2 (its series already
# ended before this column, so it stays blank) and row 74, whose figure
# was revised with this release.
$ws.Range("BH3:BH73").Value2 = $ws.Range("BG3:BG73").Value2

# Row 74 (01-01-2021) was revised with this release.
$ws.Range("BH74").Value2 = 39368

# --- New row 75 for the freshly published quarter 01-04-2021 -----------
# Write the date label as plain text (matching the style of the other
# "Serie" labels in column A, which carry no explicit number format), not
# as an auto-recognised date value.
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value2 = "01-04-2021"
$ws.Range("A75").ClearFormats()

$ws.Range("BH75").Value2 = 41667

$excel.CutCopyMode = 0
